$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KEYS")

# Row 4 (Inventory)
$ws.Range("C4").Value = 757000000.0
$ws.Range("D4").Value = 778000000.0
$ws.Range("E4").Value = 756000000.0
$ws.Range("F4").Value = 737000000.0

# Row 14 (Accounts Payable)
$ws.Range("B14").Value = 228000000.0
$ws.Range("C14").Value = 224000000.0
$ws.Range("D14").Value = 193000000.0
$ws.Range("E14").Value = 182000000.0
$ws.Range("F14").Value = 224000000.0

# Row 25 (Long Term Tax Liability (Deferred))
$ws.Range("B25").Value = -700000000.0
$ws.Range("C25").Value = -740000000.0
$ws.Range("D25").Value = -746000000.0
$ws.Range("E25").Value = -726000000.0
$ws.Range("F25").Value = -737000000.0
